$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "27.672.04"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.639.85"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "212.89"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  +0.04%  "
Set-TextValue "D8" "23.14"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.871.94"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.641.31"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +0.32%  "
Set-TextValue "D15" "0.562"
$ws.Range("E15").Value = "  -1.09%  "
Set-TextValue "D16" "64.64"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "27.649.32"
$ws.Range("E17").Value = "  +1.01%  "
Set-TextValue "D18" "230.50"
$ws.Range("E18").Value = "  +0.54%  "
Set-TextValue "D19" "7.73"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("E24").Value = "  -2.74%  "
Set-TextValue "D25" "149.72"
$ws.Range("E25").Value = "  +2.02%  "
Set-TextValue "D26" "6.94"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  -1.20%  "
Set-TextValue "D28" "15.65"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").Value = "1.444.42"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +0.49%  "
Set-TextValue "D37" "0.567"
$ws.Range("E37").Value = "  +0.79%  "
Set-TextValue "D38" "0.878"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +0.44%  "
Set-TextValue "D40" "0.901"
$ws.Range("E40").Value = "  +12.60%  "
Set-TextValue "D41" "70.33"
$ws.Range("E41").Value = "  +9.18%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "1.781.67"
$ws.Range("E47").Value = "  +0.05%  "
Set-TextValue "D48" "1.72"
$ws.Range("E48").Value = "  +3.74%  "
Set-TextValue "D49" "86.10"
$ws.Range("E50").Value = "  -0.67%  "
Set-TextValue "D51" "0.0989"
$ws.Range("E51").Value = "  +0.18%  "
